# Apply updated odds values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "S2"  = 1.5
    "T2"  = 2.37
    "S3"  = 1.47
    "J6"  = 2.4
    "N6"  = 7.5
    "Z6"  = 13
    "AH6" = 23
    "AI6" = 17
    "AK6" = 41
    "AO6" = 9.5
    "AW6" = 6.5
    "AX6" = 29
    "M7"  = 1.06
    "N7"  = 10
    "U7"  = 2.2
    "V7"  = 1.62
    "AG7" = 17
    "AI7" = 23
    "BA7" = 201
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
